$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fill in the new "Deposito" section (rows 16-19) on Hoja1 ---
$ws1.Range("A16").Value = "Deposito"
$ws1.Range("A17").Value = "Precio"
$ws1.Range("A18").Value = "Extras"
$ws1.Range("B18").Value = "Maceta"
$ws1.Range("C18").Value = "Plantas"
$ws1.Range("A19").Value = "Precio"
$ws1.Range("B19").Value = 495
$ws1.Range("C19").Value = 200

# Carry the green/bordered label style (used by A12/A13/.../A18) down to A19
$ws1.Range("A18").Copy()
$ws1.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Carry the currency number format (used by B15/E15/F15/...) onto B19:C19
$ws1.Range("B15").Copy()
$ws1.Range("B19:C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add the new (blank) "Hoja2" sheet after "Hoja1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"

# Keep Hoja1 as the active/selected sheet, with the same selection the
# author ended up with after entering the new data.
$ws1.Activate()
$ws1.Range("C19").Select()
